# Update Phoenix_Profits market-price derived figures across all sheets
# (scheduled runner refresh of currentAveragePrice / LevePrice* / LeveProfit* columns)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 29
$ws.Range("H29").Value = 1043.8
$ws.Range("J29").Value = 3499
$ws.Range("L29").Value = 10497
$ws.Range("N29").Value = -11059
# Row 43
$ws.Range("H43").Value = 2563.8333
$ws.Range("I43").Value = 2176.6667
$ws.Range("J43").Value = 2951
$ws.Range("K43").Value = 2176.6667
$ws.Range("L43").Value = 2951
$ws.Range("M43").Value = -2107.6667
$ws.Range("N43").Value = -3089
# Row 58
$ws.Range("H58").Value = 1622.2858
$ws.Range("I58").Value = 584.75
$ws.Range("J58").Value = 3005.6667
$ws.Range("K58").Value = 1754.25
$ws.Range("L58").Value = 9017.000100000001
$ws.Range("M58").Value = -1604.25
$ws.Range("N58").Value = -9317.000100000001
# Row 61
$ws.Range("H61").Value = 113.166664
$ws.Range("I61").Value = 113.166664
$ws.Range("K61").Value = 339.499992
$ws.Range("M61").Value = -167.499992
# Row 100
$ws.Range("H100").Value = 3243.8333
$ws.Range("I100").Value = 3158.3333
$ws.Range("K100").Value = 3158.3333
$ws.Range("M100").Value = -2617.3333
# Row 113
$ws.Range("H113").Value = 1612.875
$ws.Range("J113").Value = 1500
$ws.Range("L113").Value = 1500
$ws.Range("N113").Value = -8008

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 3590.6282
$ws.Range("I32").Value = 3280.7163
$ws.Range("K32").Value = 3280.7163
$ws.Range("M32").Value = -2993.7163
# Row 45
$ws.Range("H45").Value = 1740.5454
$ws.Range("I45").Value = 1278.6842
$ws.Range("K45").Value = 1278.6842
$ws.Range("M45").Value = -901.6841999999999
# Row 74
$ws.Range("H74").Value = 359522.7
$ws.Range("I74").Value = 520620.06
$ws.Range("K74").Value = 520620.06
$ws.Range("M74").Value = -519746.06
# Row 77
$ws.Range("H77").Value = 359522.7
$ws.Range("I77").Value = 520620.06
$ws.Range("K77").Value = 2603100.3
$ws.Range("M77").Value = -2598732.3
# Row 122
$ws.Range("H122").Value = 24774.2
$ws.Range("I122").Value = 2112.75
$ws.Range("K122").Value = 6338.25
$ws.Range("M122").Value = -3888.25
# Row 132
$ws.Range("H132").Value = 2050.9
$ws.Range("I132").Value = 1967.8718
$ws.Range("K132").Value = 5903.6154
$ws.Range("M132").Value = -3373.6154

$ws = $wb.Worksheets.Item("BSM")
# Row 107
$ws.Range("H107").Value = 2960.111
$ws.Range("I107").Value = 3436.0454
$ws.Range("K107").Value = 3436.0454
$ws.Range("M107").Value = -1516.0454

$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 6829.5186
$ws.Range("I16").Value = 5483.769
$ws.Range("J16").Value = 8079.143
$ws.Range("K16").Value = 5483.769
$ws.Range("L16").Value = 8079.143
$ws.Range("M16").Value = -5196.769
$ws.Range("N16").Value = -8653.143
# Row 99
$ws.Range("H99").Value = 4538.933
$ws.Range("I99").Value = 4652.846
$ws.Range("K99").Value = 4652.846
$ws.Range("M99").Value = -3154.846
# Row 105
$ws.Range("H105").Value = 1929.3572
$ws.Range("J105").Value = 1806.5
$ws.Range("L105").Value = 1806.5
$ws.Range("N105").Value = -5300.5
# Row 113
$ws.Range("H113").Value = 6829.5186
$ws.Range("I113").Value = 5483.769
$ws.Range("J113").Value = 8079.143
$ws.Range("K113").Value = 5483.769
$ws.Range("L113").Value = 8079.143
$ws.Range("M113").Value = -3313.769
$ws.Range("N113").Value = -12419.143
# Row 126
$ws.Range("H126").Value = 4538.933
$ws.Range("I126").Value = 4652.846
$ws.Range("K126").Value = 13958.538
$ws.Range("M126").Value = -11488.538

$ws = $wb.Worksheets.Item("CUL")
# Row 119
$ws.Range("H119").Value = 5446.5454
$ws.Range("I119").Value = 5446.5454
$ws.Range("K119").Value = 16339.6362
$ws.Range("M119").Value = -11501.6362
# Row 127
$ws.Range("H127").Value = 94770.17999999999
$ws.Range("J127").Value = 94770.17999999999
$ws.Range("L127").Value = 284310.54
$ws.Range("N127").Value = -294230.54

$ws = $wb.Worksheets.Item("GSM")
# Row 18
$ws.Range("H18").Value = 16513.857
$ws.Range("J18").Value = 15697.5
$ws.Range("L18").Value = 15697.5
$ws.Range("N18").Value = -16283.5
# Row 22
$ws.Range("H22").Value = 8813
$ws.Range("I22").Value = 5021.6665
$ws.Range("J22").Value = 14500
$ws.Range("K22").Value = 5021.6665
$ws.Range("L22").Value = 14500
$ws.Range("M22").Value = -4492.6665
$ws.Range("N22").Value = -15558
# Row 86
$ws.Range("H86").Value = 40896.332
$ws.Range("J86").Value = 40896.332
$ws.Range("L86").Value = 40896.332
$ws.Range("N86").Value = -43268.332
# Row 89
$ws.Range("H89").Value = 40896.332
$ws.Range("J89").Value = 40896.332
$ws.Range("L89").Value = 122688.996
$ws.Range("N89").Value = -134544.996
# Row 123
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 20447.586
$ws.Range("I7").Value = 29101
$ws.Range("K7").Value = 29101
$ws.Range("M7").Value = -28989
# Row 16
$ws.Range("H16").Value = 1201.22
$ws.Range("I16").Value = 1149.275
$ws.Range("J16").Value = 1409
$ws.Range("K16").Value = 1149.275
$ws.Range("L16").Value = 1409
$ws.Range("M16").Value = -979.2750000000001
$ws.Range("N16").Value = -1749
# Row 42
$ws.Range("H42").Value = 0
$ws.Range("I42").Value = 0
$ws.Range("J42").Value = 0
$ws.Range("K42").Value = 0
$ws.Range("L42").Value = 0
$ws.Range("M42").ClearContents()
$ws.Range("N42").ClearContents()
# Row 49
$ws.Range("H49").Value = 0
$ws.Range("I49").Value = 0
$ws.Range("J49").Value = 0
$ws.Range("K49").Value = 0
$ws.Range("L49").Value = 0
$ws.Range("M49").ClearContents()
$ws.Range("N49").ClearContents()
# Row 122
$ws.Range("H122").Value = 3530.9285
$ws.Range("I122").Value = 3511.2083
$ws.Range("J122").Value = 3649.25
$ws.Range("K122").Value = 10533.6249
$ws.Range("L122").Value = 10947.75
$ws.Range("M122").Value = -8083.624899999999
$ws.Range("N122").Value = -15847.75
# Row 126
$ws.Range("H126").Value = 20447.586
$ws.Range("I126").Value = 29101
$ws.Range("K126").Value = 87303
$ws.Range("M126").Value = -84833
# Row 132
$ws.Range("H132").Value = 2994.5386
$ws.Range("I132").Value = 2113.4
$ws.Range("K132").Value = 6340.200000000001
$ws.Range("M132").Value = -3810.200000000001
# Row 136
$ws.Range("H136").Value = 21311.688
$ws.Range("I136").Value = 2216.0527
$ws.Range("K136").Value = 6648.158100000001
$ws.Range("M136").Value = -4098.158100000001

$ws = $wb.Worksheets.Item("WVR")
# Row 126
$ws.Range("H126").Value = 8226.154
$ws.Range("I126").Value = 8218.25
$ws.Range("K126").Value = 24654.75
$ws.Range("M126").Value = -22184.75
# Row 136
$ws.Range("H136").Value = 55100.973
$ws.Range("I136").Value = 37134.89
$ws.Range("J136").Value = 108999.22
$ws.Range("K136").Value = 111404.67
$ws.Range("L136").Value = 326997.66
$ws.Range("M136").Value = -108854.67
$ws.Range("N136").Value = -332097.66
